# "more work on skip connection and upsampling"
#
# 1. block5c_project_bn textbox grows a second line "block5c_drop" and
#    the autofit height is corrected to match.
# 2. Four new red "Right Arrow" shapes (skip-connection / upsampling
#    arrows) are added near the top of the z-order stack, after the
#    "Bottle Neck" rectangle.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1) block5c_project_bn label gains a "block5c_drop" second line.
# ---------------------------------------------------------------------
$labelShape = $s.Shapes.Item(8)          # id=16, "TextBox 15"
$cr = [char]13
$labelShape.TextFrame.TextRange.Text = $labelShape.TextFrame.TextRange.Text + $cr + "block5c_drop"
$labelShape.Height = 36.35157775878906   # -> cy = 461665 EMU

# ---------------------------------------------------------------------
# 2) Four new red right-arrow shapes, cloned from an existing arrow so
#    they pick up the same shape "style" (lnRef/fillRef/effectRef).
# ---------------------------------------------------------------------
$templateArrow = $s.Shapes.Item(10)      # id=18, "Arrow: Right 17"

$newArrows = @()
for ($i = 0; $i -lt 5; $i++) {
    $d = $templateArrow.Duplicate()
    $newArrows += $d.Item(1)
}

# The first duplicate (id=2) is only scratch space used to burn through
# the id allocator so the remaining four line up on 5/7/9/11, exactly
# like the authored deck.
$newArrows[0].Delete()

$arrowDefs = @(
    @{ Shape = $newArrows[1]; Name = "Arrow: Right 4";  Left = 295.4992980957031;  Top = 358.10528564453125;  Width = 298.500732421875; Height = 14.611260414123535 },
    @{ Shape = $newArrows[2]; Name = "Arrow: Right 6";  Left = 295.4992980957031;  Top = 271.46112060546875;  Width = 298.500732421875; Height = 14.611260414123535 },
    @{ Shape = $newArrows[3]; Name = "Arrow: Right 8";  Left = 291.8287658691406;  Top = 141.0609588623047;   Width = 298.500732421875; Height = 14.611260414123535 },
    @{ Shape = $newArrows[4]; Name = "Arrow: Right 10"; Left = 291.82867431640625; Top = 70.94969177246094;   Width = 298.500732421875; Height = 14.611260414123535 }
)

foreach ($def in $arrowDefs) {
    $shp = $def.Shape
    $shp.Name = $def.Name
    $shp.Left = $def.Left
    $shp.Top = $def.Top
    $shp.Width = $def.Width
    $shp.Height = $def.Height
    $shp.Fill.ForeColor.RGB = 255   # RGB(255,0,0) -> srgbClr FF0000
}
